# Updates cryptos list with refreshed prices / 1h volume deltas.
# Price cells in column D are stored as plain text (e.g. "43.951.75",
# "1.00") rather than numbers, since some prices use "." as a thousands
# separator and others rely on significant trailing zeros. For values
# that look numeric (e.g. "74.55", "1.00"), a leading apostrophe is used
# so Excel keeps them as text instead of silently coercing them into
# Number cells (which would strip the formatting / trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.951.75'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').Value = '2.366.87'
$ws.Range('E3').Value = '  +0.29%  '

$ws.Range('E4').Value = '  -0.24%  '

$ws.Range('E5').Value = '  -1.38%  '

$ws.Range('E6').Value = '  +0.03%  '

$ws.Range('D7').Value = '''74.55'
$ws.Range('E7').Value = '  +0.82%  '

$ws.Range('D9').Value = '''0.603'
$ws.Range('E9').Value = '  +1.11%  '

$ws.Range('E10').Value = '  +1.87%  '

$ws.Range('D11').Value = '''60.22'
$ws.Range('E11').Value = '  +5.23%  '

$ws.Range('D12').Value = '''37.38'
$ws.Range('E12').Value = '  +16.00%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.109'
$ws.Range('E13').Value = '  +0.91%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''7.29'
$ws.Range('E14').Value = '  +0.20%  '

$ws.Range('D15').Value = '''16.46'
$ws.Range('E15').Value = '  -0.69%  '

$ws.Range('E16').Value = '  +3.25%  '

$ws.Range('D17').Value = '2.363.21'
$ws.Range('E17').Value = '  -0.02%  '

$ws.Range('D18').Value = '43.844.27'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('E19').Value = '  +2.54%  '

$ws.Range('D20').Value = '''6.64'
$ws.Range('E20').Value = '  -5.36%  '

$ws.Range('D21').Value = '''77.44'
$ws.Range('E21').Value = '  +0.18%  '

$ws.Range('D22').Value = '''254.39'
$ws.Range('E22').Value = '  -1.73%  '

$ws.Range('E23').Value = '  +3.81%  '

$ws.Range('E24').Value = '  -0.01%  '

$ws.Range('E25').Value = '  -5.99%  '

$ws.Range('E26').Value = '  +0.36%  '

$ws.Range('D27').Value = '''10.63'
$ws.Range('E27').Value = '  -1.67%  '

$ws.Range('E28').Value = '  +0.87%  '

$ws.Range('D29').Value = '''22.42'
$ws.Range('E29').Value = '  -1.63%  '

$ws.Range('D30').Value = '''175.28'
$ws.Range('E30').Value = '  -0.21%  '

$ws.Range('E31').Value = '  +0.90%  '

$ws.Range('E32').Value = '  -1.68%  '

$ws.Range('D33').Value = '''0.0763'
$ws.Range('E33').Value = '  +0.18%  '

$ws.Range('D34').Value = '''5.48'
$ws.Range('E34').Value = '  -2.16%  '

$ws.Range('D35').Value = '''5.13'
$ws.Range('E35').Value = '  -2.12%  '

$ws.Range('D36').Value = '''3.82'
$ws.Range('E36').Value = '  +1.39%  '

$ws.Range('D37').Value = '''6.63'
$ws.Range('E37').Value = '  +4.12%  '

$ws.Range('E38').Value = '  +2.05%  '

$ws.Range('E39').Value = '  +0.33%  '

$ws.Range('D40').Value = '''5.64'
$ws.Range('E40').Value = '  +17.89%  '

$ws.Range('D41').Value = '''20.79'
$ws.Range('E41').Value = '  +9.42%  '

$ws.Range('D42').Value = '''65.10'
$ws.Range('E42').Value = '  +9.92%  '

$ws.Range('D43').Value = '''0.108'
$ws.Range('E43').Value = '  -4.05%  '

$ws.Range('D44').Value = '''9.10'
$ws.Range('E44').Value = '  +1.10%  '

$ws.Range('E45').Value = '  -0.96%  '

$ws.Range('D46').Value = '''2.53'
$ws.Range('E46').Value = '  +1.17%  '

$ws.Range('E47').Value = '  +0.07%  '

$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  -0.06%  '

$ws.Range('E49').Value = '  -0.78%  '

$ws.Range('D50').Value = '''98.78'
$ws.Range('E50').Value = '  -2.08%  '

$ws.Range('E51').Value = '  +2.20%  '
